$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 11648
$ws.Range("D2").Value = 15991591
$ws.Range("C4").Value = 21924
$ws.Range("D4").Value = 27899787
$ws.Range("C7").Value = 60412
$ws.Range("D7").Value = 86857428
$ws.Range("C8").Value = 79726
$ws.Range("D8").Value = 105978151
$ws.Range("C9").Value = 29032
$ws.Range("D9").Value = 40473733
$ws.Range("C10").Value = 70219
$ws.Range("D10").Value = 100424548
$ws.Range("C11").Value = 9755
$ws.Range("D11").Value = 12742826
$ws.Range("C12").Value = 3879
$ws.Range("D12").Value = 5462503
$ws.Range("C13").Value = 15004
$ws.Range("D13").Value = 20309910
$ws.Range("C14").Value = 48103
$ws.Range("D14").Value = 64643005
$ws.Range("C15").Value = 23371
$ws.Range("D15").Value = 30391103
$ws.Range("C17").Value = 45203
$ws.Range("D17").Value = 57071623
$ws.Range("C18").Value = 47958
$ws.Range("D18").Value = 64019321
$ws.Range("C19").Value = 33620
$ws.Range("D19").Value = 41281878
$ws.Range("C20").Value = 50226
$ws.Range("D20").Value = 60975176
$ws.Range("C21").Value = 3919
$ws.Range("D21").Value = 5390326
$ws.Range("C23").Value = 6099
$ws.Range("D23").Value = 7670123
$ws.Range("C26").Value = 14628
$ws.Range("D26").Value = 20950551
$ws.Range("C27").Value = 24196
$ws.Range("D27").Value = 31675153
$ws.Range("C28").Value = 3339
$ws.Range("D28").Value = 4525503
$ws.Range("C29").Value = 20198
$ws.Range("D29").Value = 28671646
$ws.Range("C30").Value = 1792
$ws.Range("D30").Value = 2263287
$ws.Range("C32").Value = 3264
$ws.Range("D32").Value = 4314588
$ws.Range("C33").Value = 9227
$ws.Range("D33").Value = 12340150
$ws.Range("C34").Value = 4715
$ws.Range("D34").Value = 6005887
$ws.Range("C36").Value = 6840
$ws.Range("D36").Value = 8103150
$ws.Range("C37").Value = 10524
$ws.Range("D37").Value = 13923514
$ws.Range("C38").Value = 8568
$ws.Range("D38").Value = 10212876
$ws.Range("C39").Value = 15149
$ws.Range("D39").Value = 18540170
$ws.Range("C40").Value = 3208
$ws.Range("D40").Value = 4383382
$ws.Range("C42").Value = 8138
$ws.Range("D42").Value = 9941975
$ws.Range("C45").Value = 16135
$ws.Range("D45").Value = 23000055
$ws.Range("C46").Value = 25766
$ws.Range("D46").Value = 34325407
$ws.Range("C47").Value = 4483
$ws.Range("D47").Value = 6186232
$ws.Range("C48").Value = 26509
$ws.Range("D48").Value = 37913164
$ws.Range("C49").Value = 3013
$ws.Range("D49").Value = 3813847
$ws.Range("C51").Value = 4870
$ws.Range("D51").Value = 6355653
$ws.Range("C52").Value = 13751
$ws.Range("D52").Value = 18396803
$ws.Range("C53").Value = 5922
$ws.Range("D53").Value = 7390764
$ws.Range("C55").Value = 7643
$ws.Range("D55").Value = 9242396
$ws.Range("C56").Value = 17443
$ws.Range("D56").Value = 23426312
$ws.Range("C57").Value = 10467
$ws.Range("D57").Value = 12537480
$ws.Range("C58").Value = 16953
$ws.Range("D58").Value = 20848001
$ws.Range("C59").Value = 2810
$ws.Range("D59").Value = 3860409
$ws.Range("C60").Value = 5152
$ws.Range("D60").Value = 6511552
$ws.Range("C63").Value = 14411
$ws.Range("D63").Value = 20582457
$ws.Range("C64").Value = 19730
$ws.Range("D64").Value = 25690890
$ws.Range("C65").Value = 5161
$ws.Range("D65").Value = 7209214
$ws.Range("C66").Value = 16325
$ws.Range("D66").Value = 23338864
$ws.Range("C67").Value = 1987
$ws.Range("D67").Value = 2584669
$ws.Range("C69").Value = 3351
$ws.Range("D69").Value = 4485741
$ws.Range("C70").Value = 8694
$ws.Range("D70").Value = 11690954
$ws.Range("C71").Value = 4797
$ws.Range("D71").Value = 6109600
$ws.Range("C73").Value = 5538
$ws.Range("D73").Value = 6807218
$ws.Range("C74").Value = 8962
$ws.Range("D74").Value = 11808552
$ws.Range("C75").Value = 7608
$ws.Range("D75").Value = 9330182
$ws.Range("C76").Value = 13571
$ws.Range("D76").Value = 16610384
$ws.Range("C77").Value = 2646
$ws.Range("D77").Value = 3657020
$ws.Range("C78").Value = 2051
$ws.Range("D78").Value = 2711724
$ws.Range("C80").Value = 4834
$ws.Range("D80").Value = 6854359
$ws.Range("C81").Value = 5372
$ws.Range("D81").Value = 7481835
$ws.Range("C82").Value = 962
$ws.Range("D82").Value = 1375354
$ws.Range("C83").Value = 6574
$ws.Range("D83").Value = 9435455
$ws.Range("C84").Value = 388
$ws.Range("D84").Value = 510975
$ws.Range("C86").Value = 1485
$ws.Range("D86").Value = 2032048
$ws.Range("C87").Value = 3644
$ws.Range("D87").Value = 5064572
$ws.Range("C88").Value = 2144
$ws.Range("D88").Value = 2754835
$ws.Range("C89").Value = 1110
$ws.Range("D89").Value = 1356249
$ws.Range("C91").Value = 1396
$ws.Range("D91").Value = 1827659
$ws.Range("C92").Value = 2803
$ws.Range("D92").Value = 3374801
$ws.Range("C93").Value = 5418
$ws.Range("D93").Value = 7548418
$ws.Range("C95").Value = 10647
$ws.Range("D95").Value = 13796889
$ws.Range("C97").Value = 1287
$ws.Range("D97").Value = 1588130
$ws.Range("C98").Value = 30375
$ws.Range("D98").Value = 43627310
$ws.Range("C99").Value = 44909
$ws.Range("D99").Value = 59426534
$ws.Range("C100").Value = 9640
$ws.Range("D100").Value = 13212373
$ws.Range("C101").Value = 33281
$ws.Range("D101").Value = 47756886
$ws.Range("C102").Value = 3934
$ws.Range("D102").Value = 5065396
$ws.Range("C103").Value = 1928
$ws.Range("D103").Value = 2669692
$ws.Range("C104").Value = 5847
$ws.Range("D104").Value = 7841411
$ws.Range("C105").Value = 20395
$ws.Range("D105").Value = 27237444
$ws.Range("C106").Value = 9155
$ws.Range("D106").Value = 11742796
$ws.Range("C108").Value = 11463
$ws.Range("D108").Value = 13826029
$ws.Range("C109").Value = 21690
$ws.Range("D109").Value = 29481438
$ws.Range("C110").Value = 14201
$ws.Range("D110").Value = 16943714
$ws.Range("C111").Value = 30943
$ws.Range("D111").Value = 37083451
$ws.Range("C113").Value = 7705
$ws.Range("D113").Value = 10558006
$ws.Range("C115").Value = 3259
$ws.Range("D115").Value = 4573787
$ws.Range("C118").Value = 5542
$ws.Range("D118").Value = 8012228
$ws.Range("C119").Value = 8883
$ws.Range("D119").Value = 12340574
$ws.Range("C120").Value = 2404
$ws.Range("D120").Value = 3365991
$ws.Range("C121").Value = 7308
$ws.Range("D121").Value = 10441001
$ws.Range("C122").Value = 816
$ws.Range("D122").Value = 1125148
$ws.Range("C124").Value = 1224
$ws.Range("D124").Value = 1687753
$ws.Range("C125").Value = 3243
$ws.Range("D125").Value = 4552173
$ws.Range("C126").Value = 4295
$ws.Range("D126").Value = 5909762
$ws.Range("C127").Value = 1938
$ws.Range("D127").Value = 2470691
$ws.Range("C128").Value = 2324
$ws.Range("D128").Value = 3281959
$ws.Range("C129").Value = 1632
$ws.Range("D129").Value = 2214657
$ws.Range("C130").Value = 3130
$ws.Range("D130").Value = 4081274
$ws.Range("C132").Value = 3082
$ws.Range("D132").Value = 4924544
$ws.Range("C133").Value = 67
$ws.Range("D133").Value = 130058
$ws.Range("C134").Value = 1412
$ws.Range("D134").Value = 2486756
$ws.Range("C136").Value = 3813
$ws.Range("D136").Value = 7216837
$ws.Range("C137").Value = 3360
$ws.Range("D137").Value = 5985921
$ws.Range("C138").Value = 1129
$ws.Range("D138").Value = 2112364
$ws.Range("C139").Value = 2378
$ws.Range("D139").Value = 4334171
$ws.Range("C140").Value = 300
$ws.Range("D140").Value = 542569
$ws.Range("C142").Value = 283
$ws.Range("D142").Value = 492934
$ws.Range("C143").Value = 1385
$ws.Range("D143").Value = 2418235
$ws.Range("C144").Value = 1485
$ws.Range("D144").Value = 2675991
$ws.Range("C145").Value = 586
$ws.Range("D145").Value = 963223
$ws.Range("C146").Value = 636
$ws.Range("D146").Value = 1055955
$ws.Range("C147").Value = 557
$ws.Range("D147").Value = 976052
$ws.Range("C148").Value = 1078
$ws.Range("D148").Value = 1795400
$ws.Range("C149").Value = 2914
$ws.Range("D149").Value = 3969380
$ws.Range("C151").Value = 8163
$ws.Range("D151").Value = 10087176
$ws.Range("C154").Value = 25000
$ws.Range("D154").Value = 35515926
$ws.Range("C155").Value = 42868
$ws.Range("D155").Value = 54425467
$ws.Range("C156").Value = 15321
$ws.Range("D156").Value = 21173804
$ws.Range("C157").Value = 32126
$ws.Range("D157").Value = 46124555
$ws.Range("C158").Value = 3951
$ws.Range("D158").Value = 5104940
$ws.Range("C159").Value = 1925
$ws.Range("D159").Value = 2698359
$ws.Range("C160").Value = 5279
$ws.Range("D160").Value = 7070300
$ws.Range("C161").Value = 19636
$ws.Range("D161").Value = 26553545
$ws.Range("C162").Value = 9006
$ws.Range("D162").Value = 11356259
$ws.Range("C164").Value = 10057
$ws.Range("D164").Value = 12328349
$ws.Range("C165").Value = 22245
$ws.Range("D165").Value = 30092024
$ws.Range("C166").Value = 13942
$ws.Range("D166").Value = 16894870
$ws.Range("C167").Value = 29360
$ws.Range("D167").Value = 34409086
$ws.Range("C168").Value = 889
$ws.Range("D168").Value = 1226039
$ws.Range("C170").Value = 22604
$ws.Range("D170").Value = 30111193
$ws.Range("C171").Value = 417
$ws.Range("D171").Value = 608866
$ws.Range("C172").Value = 1010
$ws.Range("D172").Value = 1380815
$ws.Range("C173").Value = 66846
$ws.Range("D173").Value = 96393868
$ws.Range("C174").Value = 122570
$ws.Range("D174").Value = 167854531
$ws.Range("C175").Value = 218953
$ws.Range("D175").Value = 313564072
$ws.Range("C176").Value = 90191
$ws.Range("D176").Value = 132296241
$ws.Range("C177").Value = 41512
$ws.Range("D177").Value = 56790652
$ws.Range("C178").Value = 8737
$ws.Range("D178").Value = 12503441
$ws.Range("C179").Value = 21889
$ws.Range("D179").Value = 30916185
$ws.Range("C180").Value = 144552
$ws.Range("D180").Value = 198277826
$ws.Range("C181").Value = 43074
$ws.Range("D181").Value = 58017200
$ws.Range("C183").Value = 41402
$ws.Range("D183").Value = 50646803
$ws.Range("C184").Value = 58862
$ws.Range("D184").Value = 79318714
$ws.Range("C185").Value = 75187
$ws.Range("D185").Value = 96528207
$ws.Range("C186").Value = 63878
$ws.Range("D186").Value = 83074521
$ws.Range("C187").Value = 5160
$ws.Range("D187").Value = 6850222
$ws.Range("C189").Value = 4191
$ws.Range("D189").Value = 5595908
$ws.Range("C192").Value = 7657
$ws.Range("D192").Value = 11124470
$ws.Range("C193").Value = 13781
$ws.Range("D193").Value = 18709585
$ws.Range("C194").Value = 1995
$ws.Range("D194").Value = 2784772
$ws.Range("C195").Value = 8278
$ws.Range("D195").Value = 11694891
$ws.Range("C196").Value = 1026
$ws.Range("D196").Value = 1386088
$ws.Range("C197").Value = 408
$ws.Range("D197").Value = 589998
$ws.Range("C198").Value = 1691
$ws.Range("D198").Value = 2336764
$ws.Range("C199").Value = 4416
$ws.Range("D199").Value = 6176886
$ws.Range("C200").Value = 2787
$ws.Range("D200").Value = 3839967
$ws.Range("C201").Value = 3708
$ws.Range("D201").Value = 4796957
$ws.Range("C202").Value = 5242
$ws.Range("D202").Value = 7404077
$ws.Range("C203").Value = 2944
$ws.Range("D203").Value = 3881460
$ws.Range("C204").Value = 5192
$ws.Range("D204").Value = 6603702
$ws.Range("C205").Value = 1952
$ws.Range("D205").Value = 2481708
$ws.Range("C206").Value = 2092
$ws.Range("D206").Value = 2814768
$ws.Range("C209").Value = 3191
$ws.Range("D209").Value = 4607382
$ws.Range("C210").Value = 5582
$ws.Range("D210").Value = 7662348
$ws.Range("C211").Value = 2051
$ws.Range("D211").Value = 2882483
$ws.Range("C212").Value = 3715
$ws.Range("D212").Value = 5324075
$ws.Range("C213").Value = 571
$ws.Range("D213").Value = 770899
$ws.Range("C215").Value = 721
$ws.Range("D215").Value = 1012263
$ws.Range("C216").Value = 2600
$ws.Range("D216").Value = 3587392
$ws.Range("C217").Value = 2682
$ws.Range("D217").Value = 3672928
$ws.Range("C218").Value = 1394
$ws.Range("D218").Value = 1807220
$ws.Range("C219").Value = 1941
$ws.Range("D219").Value = 2711402
$ws.Range("C220").Value = 1185
$ws.Range("D220").Value = 1610751
$ws.Range("C221").Value = 2820
$ws.Range("D221").Value = 3685759
$ws.Range("C223").Value = 4526
$ws.Range("D223").Value = 8535478
$ws.Range("C224").Value = 853
$ws.Range("D224").Value = 1583356
$ws.Range("C226").Value = 2796
$ws.Range("D226").Value = 5302707
$ws.Range("C227").Value = 12663
$ws.Range("D227").Value = 23319110
$ws.Range("C228").Value = 2374
$ws.Range("D228").Value = 4389369
$ws.Range("C229").Value = 1359
$ws.Range("D229").Value = 2522661
$ws.Range("C233").Value = 459
$ws.Range("D233").Value = 806729
$ws.Range("C234").Value = 483
$ws.Range("D234").Value = 896477
$ws.Range("C235").Value = 330
$ws.Range("D235").Value = 610841
$ws.Range("C236").Value = 244
$ws.Range("D236").Value = 414430
$ws.Range("C237").Value = 218
$ws.Range("D237").Value = 394624
$ws.Range("C238").Value = 446
$ws.Range("D238").Value = 786551
$ws.Range("C239").Value = 3286
$ws.Range("D239").Value = 4476989
$ws.Range("C241").Value = 5697
$ws.Range("D241").Value = 7189577
$ws.Range("C244").Value = 15544
$ws.Range("D244").Value = 22225902
$ws.Range("C245").Value = 26923
$ws.Range("D245").Value = 35051551
$ws.Range("C246").Value = 5373
$ws.Range("D246").Value = 7429135
$ws.Range("C247").Value = 21835
$ws.Range("D247").Value = 31227935
$ws.Range("C248").Value = 2137
$ws.Range("D248").Value = 2693221
$ws.Range("C249").Value = 1199
$ws.Range("D249").Value = 1667407
$ws.Range("C250").Value = 3943
$ws.Range("D250").Value = 5210372
$ws.Range("C251").Value = 11588
$ws.Range("D251").Value = 15611823
$ws.Range("C252").Value = 5566
$ws.Range("D252").Value = 7061125
$ws.Range("C254").Value = 6553
$ws.Range("D254").Value = 7791529
$ws.Range("C255").Value = 10628
$ws.Range("D255").Value = 14006969
$ws.Range("C256").Value = 9075
$ws.Range("D256").Value = 11044975
$ws.Range("C257").Value = 17921
$ws.Range("D257").Value = 21717146
$ws.Range("C258").Value = 12451
$ws.Range("D258").Value = 17218437
$ws.Range("C260").Value = 17445
$ws.Range("D260").Value = 21404178
$ws.Range("C262").Value = 941
$ws.Range("D262").Value = 1070730
$ws.Range("C263").Value = 45072
$ws.Range("D263").Value = 63914862
$ws.Range("C264").Value = 66138
$ws.Range("D264").Value = 86700093
$ws.Range("C265").Value = 13664
$ws.Range("D265").Value = 18826318
$ws.Range("C266").Value = 48690
$ws.Range("D266").Value = 68692136
$ws.Range("C267").Value = 6221
$ws.Range("D267").Value = 8045311
$ws.Range("C268").Value = 2937
$ws.Range("D268").Value = 4108131
$ws.Range("C269").Value = 11486
$ws.Range("D269").Value = 15160653
$ws.Range("C270").Value = 33528
$ws.Range("D270").Value = 44885774
$ws.Range("C271").Value = 17377
$ws.Range("D271").Value = 21919642
$ws.Range("C273").Value = 19160
$ws.Range("D273").Value = 22258487
$ws.Range("C274").Value = 34562
$ws.Range("D274").Value = 45519594
$ws.Range("C275").Value = 22508
$ws.Range("D275").Value = 27082417
$ws.Range("C276").Value = 41318
$ws.Range("D276").Value = 49321150
$ws.Range("C277").Value = 13500
$ws.Range("D277").Value = 18177995
$ws.Range("C279").Value = 18992
$ws.Range("D279").Value = 23479161
$ws.Range("C282").Value = 58091
$ws.Range("D282").Value = 82281833
$ws.Range("C283").Value = 74786
$ws.Range("D283").Value = 98071829
$ws.Range("C284").Value = 15827
$ws.Range("D284").Value = 21472556
$ws.Range("C285").Value = 59057
$ws.Range("D285").Value = 83518799
$ws.Range("C286").Value = 7452
$ws.Range("D286").Value = 9513933
$ws.Range("C287").Value = 3088
$ws.Range("D287").Value = 4268743
$ws.Range("C288").Value = 13048
$ws.Range("D288").Value = 17385753
$ws.Range("C289").Value = 38482
$ws.Range("D289").Value = 51888426
$ws.Range("C290").Value = 18918
$ws.Range("D290").Value = 23903923
$ws.Range("C292").Value = 23566
$ws.Range("D292").Value = 27534527
$ws.Range("C293").Value = 37883
$ws.Range("D293").Value = 50134602
$ws.Range("C294").Value = 26659
$ws.Range("D294").Value = 32056520
$ws.Range("C295").Value = 42580
$ws.Range("D295").Value = 49455505
$ws.Range("C296").Value = 4084
$ws.Range("D296").Value = 5703439
$ws.Range("C298").Value = 7834
$ws.Range("D298").Value = 9659436
$ws.Range("C301").Value = 18175
$ws.Range("D301").Value = 26031995
$ws.Range("C302").Value = 28735
$ws.Range("D302").Value = 37525856
$ws.Range("C303").Value = 7060
$ws.Range("D303").Value = 9867777
$ws.Range("C304").Value = 23474
$ws.Range("D304").Value = 33577805
$ws.Range("C305").Value = 3542
$ws.Range("D305").Value = 4554422
$ws.Range("C307").Value = 5662
$ws.Range("D307").Value = 7561869
$ws.Range("C308").Value = 17752
$ws.Range("D308").Value = 23905422
$ws.Range("C309").Value = 6742
$ws.Range("D309").Value = 8730398
$ws.Range("C310").Value = 8678
$ws.Range("D310").Value = 10349134
$ws.Range("C311").Value = 19259
$ws.Range("D311").Value = 25335303
$ws.Range("C312").Value = 12469
$ws.Range("D312").Value = 15217322
$ws.Range("C313").Value = 20736
$ws.Range("D313").Value = 25425147
$ws.Range("C314").Value = 5558
$ws.Range("D314").Value = 7550272
$ws.Range("C316").Value = 18103
$ws.Range("D316").Value = 23387065
$ws.Range("C319").Value = 50973
$ws.Range("D319").Value = 72658607
$ws.Range("C320").Value = 79077
$ws.Range("D320").Value = 105020455
$ws.Range("C321").Value = 27340
$ws.Range("D321").Value = 38398774
$ws.Range("C322").Value = 56317
$ws.Range("D322").Value = 81040343
$ws.Range("C323").Value = 8049
$ws.Range("D323").Value = 10514965
$ws.Range("C324").Value = 3477
$ws.Range("D324").Value = 4893212
$ws.Range("C325").Value = 16194
$ws.Range("D325").Value = 22306476
$ws.Range("C326").Value = 41378
$ws.Range("D326").Value = 56017482
$ws.Range("C327").Value = 23335
$ws.Range("D327").Value = 30332887
$ws.Range("C329").Value = 25045
$ws.Range("D329").Value = 30209055
$ws.Range("C330").Value = 34794
$ws.Range("D330").Value = 46194345
$ws.Range("C331").Value = 24796
$ws.Range("D331").Value = 30885439
$ws.Range("C332").Value = 42625
$ws.Range("D332").Value = 50797234
